# Implemented new way of specifying flags. Added test cases.
#
# - Uppercase the language / mode flags in the "Sheet2" header row
#   ("en"/"pl"/"ua" -> "EN"/"PL"/"UA", "import" -> "IMPORT").
# - Give columns D:F their own (wider) widths instead of sharing one.
# - Move the remembered cell-selection on Sheet1 and Sheet2.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Update header labels to the new upper-case flag format -------------
$ws2.Range("A1").Value = "phrase (EN)"
$ws2.Range("B1").Value = "phrase (PL)"
$ws2.Range("C1").Value = "phrase (UA)"
$ws2.Range("D1").Value = "picture (EN, IMPORT)"
$ws2.Range("E1").Value = "picture (PL, IMPORT)"
$ws2.Range("F1").Value = "picture (UA, IMPORT)"

# --- Give columns D, E and F their own explicit widths -------------------
$ws2.Columns("D").ColumnWidth = 19.91
$ws2.Columns("E").ColumnWidth = 19.77
$ws2.Columns("F").ColumnWidth = 19.91

# --- Move the active-cell selections -------------------------------------
# (select Sheet2 first so that Sheet1 ends up as the active/selected tab,
#  matching the original tabSelected state)
$ws2.Range("E2").Select() | Out-Null
$ws1.Range("C1").Select() | Out-Null
